$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the shared-string literal in A1 with the TRUE() formula.
$ws.Range("A1").Formula = "=TRUE()"

# Move/record the active selection on Sheet1 to A1 (matches the saved view state).
$ws.Range("A1").Select() | Out-Null
